$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column C (iconResource), shifting
# iconResource/title/price/etc right by one. This makes room for the
# new "type" column.
$ws.Columns("C").Insert()

# Header row
$ws.Range("C1").Value = "type"

# New "type" column values mirror the "level" column (B) for each data row.
$ws.Range("C2").Value = 0
$ws.Range("C3").Value = 1
$ws.Range("C4").Value = 2
$ws.Range("C5").Value = 0
$ws.Range("C6").Value = 1
$ws.Range("C7").Value = 2

# Rows 4 and 7 use the wrap-text style (matching column D in those rows);
# the column insert above only inherited formatting from column B.
$ws.Range("C4").WrapText = $true
$ws.Range("C7").WrapText = $true

# New column C should carry the same custom width as column D (the old
# iconResource column it was inserted in front of, which keeps its
# original ~17.93-character width after the insert shifts it over).
$ws.Columns("C").ColumnWidth = 17.22

# Move the selection to match the recorded post-edit state.
$ws.Range("D8").Select()
